$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Helper: append `$text` to the end of paragraph `$idx` (before its pilcrow)
# as its OWN run (even if neighbouring formatting is identical) by
# toggling Bold on/off right after insertion, then stamp the language back
# onto it so the run keeps the nl-NL rPr instead of an empty one.
# -----------------------------------------------------------------------

# ============================================================
# 1) Paragraph 3 ("Talententest (blz. 1 dictaat P4P)") becomes
#    two runs "Functioneel ontwerp" + "."  (the _GoBack bookmark is
#    relocated here later, once the paragraph deletions below are done
#    and the final text of this paragraph is settled).
# ============================================================
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("Talententest (blz. 1 dictaat P4P)", $true, $false, $false, $false, $false, $true, 1, $false, "Functioneel ontwerp", 2) | Out-Null

$p3 = $d.Paragraphs(3).Range
$pos = $p3.End - 1
$dotRun = $d.Range($pos, $pos)
$dotRun.InsertAfter(".")
$dotRun.LanguageID = "nl-NL"
$dotRun.Font.Bold = $true
$dotRun.Font.Bold = $false

# ============================================================
# 2) Paragraph 4 ("Functioneel ontwerp") ->
#    "Procesverslag (blz. 1 Planningsschema)"
# ============================================================
$p4 = $d.Paragraphs(4).Range
$p4.Find.Execute("Functioneel ontwerp", $true, $false, $false, $false, $false, $true, 1, $false, "Procesverslag (blz. 1 Planningsschema)", 2) | Out-Null

# ============================================================
# 3) Paragraph 5 ("Standaardverslag (blz. 13 dictaat P4P)") becomes
#    three runs: "Individueel " + "reflectie" + "verslag (blz. 1 dictaat P4P)"
# ============================================================
$p5 = $d.Paragraphs(5).Range
$p5.Find.Execute("Standaardverslag (blz. 13 dictaat P4P)", $true, $false, $false, $false, $false, $true, 1, $false, "Individueel ", 2) | Out-Null

$p5 = $d.Paragraphs(5).Range
$pos = $p5.End - 1
$run2 = $d.Range($pos, $pos)
$run2.InsertAfter("reflectie")
$run2.LanguageID = "nl-NL"
$run2.Font.Bold = $true
$run2.Font.Bold = $false

$p5 = $d.Paragraphs(5).Range
$pos = $p5.End - 1
$run3 = $d.Range($pos, $pos)
$run3.InsertAfter("verslag (blz. 1 dictaat P4P)")
$run3.LanguageID = "nl-NL"
$run3.Font.Bold = $true
$run3.Font.Bold = $false

# ============================================================
# 4) Paragraph 6 ("Meetverslag (blz. 13 dictaat P4P)") ->
#    "Groepsplanning (blz. 5 dictaat P4P)"
# ============================================================
$p6 = $d.Paragraphs(6).Range
$p6.Find.Execute("Meetverslag (blz. 13 dictaat P4P)", $true, $false, $false, $false, $false, $true, 1, $false, "Groepsplanning (blz. 5 dictaat P4P)", 2) | Out-Null

# ============================================================
# 5) Paragraph 7 ("Procesverslag (blz. 1 Planningsschema)") ->
#    "Presentatie (blz. 9 dictaat P4P)"
# ============================================================
$p7 = $d.Paragraphs(7).Range
$p7.Find.Execute("Procesverslag (blz. 1 Planningsschema)", $true, $false, $false, $false, $false, $true, 1, $false, "Presentatie (blz. 9 dictaat P4P)", 2) | Out-Null

# ============================================================
# 6) Remove the now-duplicated / obsolete list paragraphs 8-11
#    ("Individueel reflectieverslag ...", "Groepsplanning ...",
#    "Individuele planning ...", "Presentatie ...") - delete from the
#    bottom up so earlier indices stay valid.
# ============================================================
$d.Paragraphs(11).Range.Delete()
$d.Paragraphs(10).Range.Delete()
$d.Paragraphs(9).Range.Delete()
$d.Paragraphs(8).Range.Delete()

# ============================================================
# 7) Drop the old bookmark that used to sit at the very end of the
#    "Vragen voor tutor:" heading paragraph.
# ============================================================
$d.Bookmarks("_GoBack").Delete()

# Re-add it at its new home (end of paragraph 3, after "Functioneel ontwerp.")
$p3 = $d.Paragraphs(3).Range
$pos = $p3.End - 1
$placeholder = $d.Range($pos, $pos)
$placeholder.InsertAfter("X")

$p3 = $d.Paragraphs(3).Range
$bmPos = $p3.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$p3 = $d.Paragraphs(3).Range
$placeholderRange = $d.Range($p3.End - 2, $p3.End - 1)
$placeholderRange.Delete()

# ============================================================
# 8) Append a new, empty (bold, nl-NL) paragraph at the very end of the
#    document, after "16 juni 2016 presentatie & afsluiting".
# ============================================================
$lastPara = $d.Paragraphs($d.Paragraphs.Count).Range
$endOfDoc = $lastPara.End
$insertion = $d.Range($endOfDoc, $endOfDoc)
$insertion.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count).Range
$newPara.Font.Bold = $true
$newPara.LanguageID = "nl-NL"
